# Populate cell values in the "Valores" sheet automatically.
# This mirrors running the data-extraction script and writing its
# results into the worksheet, then refreshing the totals row (row 16)
# as the sum of rows 1-15 for each column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row/column (only the cells that actually changed are
# listed here; everything else in the sheet stays 0 as before).
$values = @{
    1  = @{ A=2; B=12;          E=6;  F=3;      H=5;  I=1;          L=1 }
    2  = @{          D=0; E=10; F=6; G=3; H=0;       J=7; K=1; L=3 }
    4  = @{     B=2;               F=2;      H=1 }
    5  = @{     B=1;                         H=1 }
    6  = @{     B=5;               E=6; F=1;      H=3;                L=1 }
    9  = @{              D=0;      F=1;      H=0;          J=1;      L=0 }
    10 = @{              D=0; E=1; F=0;      H=0;               K=1; L=0 }
    11 = @{                   E=1; F=1;                     J=0;      L=1 }
    13 = @{                                                  J=0 }
    15 = @{                   E=5; F=4; G=2;                 J=3;      L=1 }
    16 = @{ A=2; B=20;        D=0; E=29; F=18; G=5; H=10; I=1; J=13; K=2; L=7 }
}

foreach ($rowNum in $values.Keys) {
    $rowValues = $values[$rowNum]
    foreach ($col in $rowValues.Keys) {
        $cellRef = "$col$rowNum"
        $ws.Range($cellRef).Value = $rowValues[$col]
    }
}

$wb.Save()
